# Fill in day 1-3 (columns B, C, D) sleep-diary data for the last
# week-block ("第七天 日期:2025-12-19" section, rows 59-76) and move the
# view/selection down to where the new data was entered.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 63: 您今天早上几点醒来? (wake-up time) ---
$ws.Range("B63").Value = 0.28472222222222221
$ws.Range("B63").NumberFormat = "h:mm"
$ws.Range("C63").Value = 0.39583333333333331
$ws.Range("C63").NumberFormat = "h:mm"
$ws.Range("D63").Value = "8：14"

# --- Row 64: 您今天几点起床? (out-of-bed time) ---
$ws.Range("B64").Value = 0.28472222222222221
$ws.Range("B64").NumberFormat = "h:mm"
$ws.Range("C64").Value = 0.4375
$ws.Range("C64").NumberFormat = "h:mm"
$ws.Range("D64").Value = "8：20"

# --- Row 65: 您昨晚几点上床? (to-bed time) ---
$ws.Range("B65").Value = 0.95833333333333337
$ws.Range("B65").NumberFormat = "h:mm"
$ws.Range("C65").Value = 0.10416666666666667
$ws.Range("C65").NumberFormat = "h:mm"
$ws.Range("D65").Value = "22：30"

# --- Row 66: 您昨晚几点熄灯? (lights-off time) ---
$ws.Range("B66").Value = 0.97916666666666663
$ws.Range("B66").NumberFormat = "h:mm"
$ws.Range("C66").Value = 0.10416666666666667
$ws.Range("C66").NumberFormat = "h:mm"
$ws.Range("D66").Value = "na"

# --- Row 67: 熄灯后多久入睡 (minutes to fall asleep) ---
$ws.Range("B67").Value = 30
$ws.Range("C67").Value = 3
$ws.Range("D67").Value = "上床后30 min"

# --- Row 68: 整晚醒来几次 (times woken) ---
$ws.Range("B68").Value = 0
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 2

# --- Row 69: 整晚共醒了多长时间 (minutes awake) ---
$ws.Range("B69").Value = 0
$ws.Range("C69").Value = 0
$ws.Range("D69").Value = 65

# --- Row 70: 整晚共睡了多长时间 (minutes slept) ---
$ws.Range("B70").Value = 410
$ws.Range("C70").Value = 420
$ws.Range("D70").Value = 480

# --- Row 71: 睡前是否使用了影响睡眠的物质 ---
$ws.Range("B71").Value = "无"
$ws.Range("C71").Value = "无"
$ws.Range("D71").Value = "无"

# --- Row 72: 睡前是否使用了电子产品/多长时间 ---
$ws.Range("B72").Value = 60
$ws.Range("C72").Value = 2
$ws.Range("D72").Value = 30

# --- Row 73: 睡前身体紧张程度 ---
$ws.Range("B73").Value = 3
$ws.Range("C73").Value = 2
$ws.Range("D73").Value = 2

# --- Row 74: 睡前精神紧张程度 ---
$ws.Range("B74").Value = 3
$ws.Range("C74").Value = 2
$ws.Range("D74").Value = 4

# --- Row 75: 整晚睡眠质量 ---
$ws.Range("B75").Value = 2
$ws.Range("C75").Value = 1
$ws.Range("D75").Value = 1

# --- Row 76: 白天是否小睡/多长时间 ---
$ws.Range("B76").Value = "无"
$ws.Range("C76").Value = "无"
$ws.Range("D76").Value = "无"

# Move the viewport down to the newly filled-in area and select D76,
# matching where the user ended up after finishing data entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 84
$win.ScrollColumn = 1
$ws.Range("D76").Select()
